$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage on the Price column so numeric-looking values
# (e.g. "1.00", "0.999", "18.02") are kept as literal text instead of
# being auto-converted to numbers (which would strip formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$sub3 = [char]0x2083

$ws.Range("D2").Value = '67.516.04'
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").Value = '3.369.57'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '589.63'
$ws.Range("E5").Value = '  +6.17%  '
$ws.Range("D6").Value = '187.57'
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("D9").Value = '0.183'
$ws.Range("E9").Value = '  +1.89%  '
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").Value = '47.47'
$ws.Range("E11").Value = '  +2.47%  '
$ws.Range("D12").Value = '0.0000275'
$ws.Range("E12").Value = '  +2.18%  '
$ws.Range("D13").Value = '3.915.02'
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("B14").Value = 'BitcoinCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D14").Value = '635.79'
$ws.Range("E14").Value = '  +7.04%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '8.61'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '67.533.18'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '3.378.16'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").Value = '18.02'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").Value = '11.13'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").Value = '18.01'
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").Value = '99.41'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("E26").Value = '  +4.41%  '
$ws.Range("D27").Value = '9.71'
$ws.Range("D28").Value = '32.57'
$ws.Range("E28").Value = '  +7.14%  '
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("D30").Value = '6.94'
$ws.Range("E30").Value = '  +3.81%  '
$ws.Range("D31").Value = '607.58'
$ws.Range("E31").Value = '  +4.14%  '
$ws.Range("D32").Value = '3.81'
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").Value = '3.982.67'
$ws.Range("E33").Value = '  +7.57%  '
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("E35").Value = '  +1.90%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '56.04'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  +7.23%  '
$ws.Range("E39").Value = '  +4.46%  '
$ws.Range("D40").Value = '33.78'
$ws.Range("E40").Value = '  +0.58%  '
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("D42").Value = "0.0${sub3}0704"
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").Value = '3.40'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").Value = '0.342'
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '0.0423'
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("E46").Value = '  +1.23%  '
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("E49").Value = '  +10.76%  '
$ws.Range("D50").Value = '2.81'
$ws.Range("E50").Value = '  -19.19%  '
$ws.Range("D51").Value = '128.43'
$ws.Range("E51").Value = '  +3.48%  '
